$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 291.7143
$ws.Range("I4").Value = 257
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 257
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -143
$ws.Range("N4").Value = -728

$ws.Range("H8").Value = 46.77778
$ws.Range("I8").Value = 40.285713
$ws.Range("K8").Value = 120.857139
$ws.Range("M8").Value = 18.142861

$ws.Range("H39").Value = 1122.0625
$ws.Range("I39").Value = 79.416664
$ws.Range("K39").Value = 238.249992
$ws.Range("M39").Value = 57.75000800000001

$ws.Range("H75").Value = 70000
$ws.Range("J75").Value = 70000
$ws.Range("L75").Value = 70000
$ws.Range("N75").Value = -71872

$ws.Range("H78").Value = 70000
$ws.Range("J78").Value = 70000
$ws.Range("L78").Value = 210000
$ws.Range("N78").Value = -219360

$ws.Range("H113").Value = 5923.4165
$ws.Range("J113").Value = 6837.6
$ws.Range("L113").Value = 6837.6
$ws.Range("N113").Value = -13345.6

$ws.Range("H116").Value = 9330
$ws.Range("I116").Value = 7995
$ws.Range("K116").Value = 7995
$ws.Range("M116").Value = -4553

$ws.Range("H131").Value = 3272.3333
$ws.Range("I131").Value = 2806.625
$ws.Range("K131").Value = 8419.875
$ws.Range("M131").Value = -3379.875

$ws.Range("H137").Value = 13831.516
$ws.Range("I137").Value = 14028.923
$ws.Range("K137").Value = 42086.769
$ws.Range("M137").Value = -39536.769

$ws.Range("H138").Value = 21743416
$ws.Range("J138").Value = 43485904
$ws.Range("L138").Value = 130457712
$ws.Range("N138").Value = -130467992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2877.077
$ws.Range("I61").Value = 2882.0454
$ws.Range("J61").Value = 2849.75
$ws.Range("K61").Value = 2882.0454
$ws.Range("L61").Value = 2849.75
$ws.Range("M61").Value = -2670.0454
$ws.Range("N61").Value = -3273.75

$ws.Range("H74").Value = 171026.14
$ws.Range("I74").Value = 199044.83
$ws.Range("K74").Value = 199044.83
$ws.Range("M74").Value = -198170.83

$ws.Range("H77").Value = 171026.14
$ws.Range("I77").Value = 199044.83
$ws.Range("K77").Value = 995224.1499999999
$ws.Range("M77").Value = -990856.1499999999

$ws.Range("H110").Value = 22597.438
$ws.Range("I110").Value = 22597.438
$ws.Range("K110").Value = 22597.438
$ws.Range("M110").Value = -20552.438

$ws.Range("H122").Value = 1753.6897
$ws.Range("I122").Value = 1469.0834
$ws.Range("J122").Value = 3119.8
$ws.Range("K122").Value = 4407.2502
$ws.Range("L122").Value = 9359.400000000001
$ws.Range("M122").Value = -1957.2502
$ws.Range("N122").Value = -14259.4

$ws.Range("H132").Value = 78151.375
$ws.Range("I132").Value = 10337.4375
$ws.Range("K132").Value = 31012.3125
$ws.Range("M132").Value = -28482.3125

$ws.Range("H136").Value = 2877.077
$ws.Range("I136").Value = 2882.0454
$ws.Range("J136").Value = 2849.75
$ws.Range("K136").Value = 8646.136200000001
$ws.Range("L136").Value = 8549.25
$ws.Range("M136").Value = -6096.136200000001
$ws.Range("N136").Value = -13649.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 415.2
$ws.Range("I11").Value = 244
$ws.Range("J11").Value = 1100
$ws.Range("K11").Value = 244
$ws.Range("L11").Value = 1100
$ws.Range("M11").Value = -104
$ws.Range("N11").Value = -1380

$ws.Range("H20").Value = 1114.2188
$ws.Range("I20").Value = 1004.2
$ws.Range("J20").Value = 1507.1428
$ws.Range("K20").Value = 1004.2
$ws.Range("L20").Value = 1507.1428
$ws.Range("M20").Value = -757.2
$ws.Range("N20").Value = -2001.1428

$ws.Range("H94").Value = 1692.9231
$ws.Range("I94").Value = 1084.9231
$ws.Range("K94").Value = 1084.9231
$ws.Range("M94").Value = -633.9231

$ws.Range("H134").Value = 1817.1765
$ws.Range("J134").Value = 9506.5
$ws.Range("L134").Value = 28519.5
$ws.Range("N134").Value = -33589.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 222.15384
$ws.Range("I7").Value = 52.166668
$ws.Range("J7").Value = 367.85715
$ws.Range("K7").Value = 52.166668
$ws.Range("L7").Value = 367.85715
$ws.Range("M7").Value = 60.833332
$ws.Range("N7").Value = -593.85715

$ws.Range("H58").Value = 17188944
$ws.Range("I58").Value = 662.2222
$ws.Range("K58").Value = 662.2222
$ws.Range("M58").Value = -459.2222

$ws.Range("H94").Value = 2346.4
$ws.Range("J94").Value = 2445.75
$ws.Range("L94").Value = 2445.75
$ws.Range("N94").Value = -3347.75

$ws.Range("H107").Value = 389
$ws.Range("I107").Value = 350.75
$ws.Range("K107").Value = 350.75
$ws.Range("M107").Value = 1569.25

$ws.Range("H122").Value = 2075.7778
$ws.Range("I122").Value = 2075.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6227.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3777.3334

$ws.Range("H136").Value = 17188944
$ws.Range("I136").Value = 662.2222
$ws.Range("K136").Value = 1986.6666
$ws.Range("M136").Value = 563.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 104.333336
$ws.Range("I14").Value = 104.333336
$ws.Range("K14").Value = 313.000008
$ws.Range("M14").Value = -140.000008

$ws.Range("H120").Value = 4699.5
$ws.Range("I120").Value = 4699.5
$ws.Range("K120").Value = 14098.5
$ws.Range("M120").Value = -9260.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 272000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 272000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 107324.95
$ws.Range("I70").Value = 125089
$ws.Range("K70").Value = 125089
$ws.Range("M70").Value = -124819

$ws.Range("H73").Value = 107324.95
$ws.Range("I73").Value = 125089
$ws.Range("K73").Value = 125089
$ws.Range("M73").Value = -124153

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1911.3334
$ws.Range("I16").Value = 1329.9
$ws.Range("K16").Value = 1329.9
$ws.Range("M16").Value = -1159.9

$ws.Range("H22").Value = 2636.182
$ws.Range("I22").Value = 1633.1666
$ws.Range("J22").Value = 3839.8
$ws.Range("K22").Value = 1633.1666
$ws.Range("L22").Value = 3839.8
$ws.Range("M22").Value = -1338.1666
$ws.Range("N22").Value = -4429.8

$ws.Range("H27").Value = 2636.182
$ws.Range("I27").Value = 1633.1666
$ws.Range("J27").Value = 3839.8
$ws.Range("K27").Value = 1633.1666
$ws.Range("L27").Value = 3839.8
$ws.Range("M27").Value = -1526.1666
$ws.Range("N27").Value = -4053.8

$ws.Range("H55").Value = 731.65717
$ws.Range("I55").Value = 522.6667
$ws.Range("J55").Value = 952.94116
$ws.Range("K55").Value = 522.6667
$ws.Range("L55").Value = 952.94116
$ws.Range("M55").Value = -349.6667
$ws.Range("N55").Value = -1298.94116

$ws.Range("H58").Value = 2250
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 2000
$ws.Range("N58").Value = -2520

$ws.Range("H61").Value = 2947
$ws.Range("J61").Value = 4101.5
$ws.Range("L61").Value = 4101.5
$ws.Range("N61").Value = -4505.5

$ws.Range("H113").Value = 2947
$ws.Range("J113").Value = 4101.5
$ws.Range("L113").Value = 4101.5
$ws.Range("N113").Value = -8441.5

$ws.Range("H132").Value = 1309.125
$ws.Range("I132").Value = 1311.8
$ws.Range("K132").Value = 3935.4
$ws.Range("M132").Value = -1405.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1148.1
$ws.Range("I113").Value = 1032.5333
$ws.Range("K113").Value = 3097.5999
$ws.Range("M113").Value = -927.5999000000002

$ws.Range("H132").Value = 1913.5416
$ws.Range("I132").Value = 927.73334
$ws.Range("J132").Value = 3556.5557
$ws.Range("K132").Value = 2783.20002
$ws.Range("L132").Value = 10669.6671
$ws.Range("M132").Value = -253.2000200000002
$ws.Range("N132").Value = -15729.6671
